# Updated cryptos list with refreshed Price / Volume(1h) figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.919.93"
$ws.Range("E2").Value = "'  +2.37%  "
$ws.Range("D3").Value = "'2.948.92"
$ws.Range("E3").Value = "'  +0.63%  "
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("D5").Value = "'591.94"
$ws.Range("E5").Value = "'  -0.34%  "
$ws.Range("D6").Value = "'147.13"
$ws.Range("E6").Value = "'  +2.48%  "
$ws.Range("E7").Value = "'  -0.09%  "
$ws.Range("D8").Value = "'2.947.38"
$ws.Range("E8").Value = "'  +0.64%  "
$ws.Range("E9").Value = "'  +1.05%  "
$ws.Range("E10").Value = "'  +1.16%  "
$ws.Range("D11").Value = "'0.148"
$ws.Range("E11").Value = "'  +5.08%  "
$ws.Range("E12").Value = "'  -0.29%  "
$ws.Range("E13").Value = "'  +4.41%  "
$ws.Range("D14").Value = "'32.42"
$ws.Range("E14").Value = "'  -2.33%  "
$ws.Range("E15").Value = "'  -1.28%  "
$ws.Range("D16").Value = "'3.437.69"
$ws.Range("E16").Value = "'  +0.65%  "
$ws.Range("D17").Value = "'62.909.40"
$ws.Range("E17").Value = "'  +2.41%  "
$ws.Range("E18").Value = "'  +0.61%  "
$ws.Range("D19").Value = "'2.948.22"
$ws.Range("E19").Value = "'  +0.56%  "
$ws.Range("D20").Value = "'437.51"
$ws.Range("E20").Value = "'  +1.01%  "
$ws.Range("E21").Value = "'  -1.09%  "
$ws.Range("E22").Value = "'  -1.24%  "
$ws.Range("E23").Value = "'  -0.79%  "
$ws.Range("E24").Value = "'  +3.73%  "
$ws.Range("D25").Value = "'80.63"
$ws.Range("E25").Value = "'  -0.95%  "
$ws.Range("E26").Value = "'  -2.31%  "
$ws.Range("E27").Value = "'  +0.47%  "
$ws.Range("E28").Value = "'  +0.02%  "
$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = "'  +1.72%  "
$ws.Range("E30").Value = "'  +5.72%  "
$ws.Range("E31").Value = "'  +0.15%  "
$ws.Range("E32").Value = "'  +12.95%  "
$ws.Range("E33").Value = "'  -1.53%  "
$ws.Range("E34").Value = "'  -0.84%  "
$ws.Range("E35").Value = "'  +0.04%  "
$ws.Range("D36").Value = "'0.990"
$ws.Range("E36").Value = "'  -1.87%  "
$ws.Range("E37").Value = "'  -0.13%  "
$ws.Range("E38").Value = "'  +0.86%  "
$ws.Range("D39").Value = "'49.58"
$ws.Range("E39").Value = "'  -0.36%  "
$ws.Range("E40").Value = "'  +1.45%  "
$ws.Range("D41").Value = "'0.117"
$ws.Range("E41").Value = "'  -4.15%  "
$ws.Range("E42").Value = "'  -0.96%  "
$ws.Range("E43").Value = "'  -0.01%  "
$ws.Range("D44").Value = "'39.45"
$ws.Range("E44").Value = "'  -5.92%  "
$ws.Range("D45").Value = "'135.73"
$ws.Range("E45").Value = "'  +1.93%  "
$ws.Range("D46").Value = "'2.682.93"
$ws.Range("E46").Value = "'  -0.42%  "
$ws.Range("E47").Value = "'  -1.94%  "
$ws.Range("D48").Value = "'357.69"
$ws.Range("E48").Value = "'  -1.60%  "
$ws.Range("E50").Value = "'  -0.69%  "
$ws.Range("D51").Value = "'22.60"
$ws.Range("E51").Value = "'  -3.91%  "
